# The "order_2" sheet's A1 header was renamed from "0_x_half" to
# "0_x_lthalf" (typo fix / completed label, per the commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("order_2")
$ws.Range("A1").Value = "0_x_lthalf"

# The author was working on (and left selected) this sheet/cell when the
# file was saved - this sheet becomes the active tab, with cell C6
# selected (moving tabSelected/focus away from the previously active
# "order_3" sheet).
$ws.Activate()
$ws.Range("C6").Select()
